# "Morada" (address) column should stop right before the postal code and
# must not include the postal code / city that used to follow it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I8").Value  = "Rua Cidade de Córdova, 2A"
$ws.Range("I9").Value  = "Rua da Ermida, N.º 64"
$ws.Range("I15").Value = "Rua Ana de Castro Osório, Nº 2 C"
$ws.Range("I21").Value = "Rua António Bessa Leite, 1468, Loja 11, Porto"
$ws.Range("I22").Value = "Rua da Reserva Botânica, Número 11, Garagem 4, Rinchoa"
$ws.Range("I23").Value = "Rua Gago Coutinho e Sacadura Cabral, Nº 21/21A/21B"
$ws.Range("I25").Value = "Rua das Passarias, Nº 251"
$ws.Range("I27").Value = "Rua Manuel Teixeira Gomes, Nº. 22, 1º., Dtº."

# Remove the trailing rows (41-59) that were appended by a later scrape run.
$ws.Rows("41:59").Delete()
